# Trade #4 closed at 2026-02-17 07:57:39 - unknown UNKNOWN +0.000%
$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1199.95
$wsSummary.Range("B4").Value = -0.05
$wsSummary.Range("B5").Value = -0.25
$wsSummary.Range("B6").Value = 4
$wsSummary.Range("B8").Value = 3
$wsSummary.Range("B9").Value = 25

# --- Strategy Status sheet (MarketMaking row) ---
$wsStrategy = $wb.Worksheets.Item("Strategy Status")
$wsStrategy.Range("C4").Value = 99.95
$wsStrategy.Range("D4").Value = 4
$wsStrategy.Range("E4").Value = -0.05
$wsStrategy.Range("F4").Value = -0.05
$wsStrategy.Range("G4").Value = 25

# --- New Trade #4 row, appended identically to "All Trades" and "MarketMaking" sheets ---
function Add-TradeRow($ws) {
    $ws.Cells.Item(5, 1).Value = 4

    # Force the Date column to be stored as plain text, matching the existing
    # rows above it (which are text, not real dates).
    $ws.Cells.Item(5, 2).NumberFormat = "@"
    $ws.Cells.Item(5, 2).Value = "2026-02-17"

    $ws.Cells.Item(5, 3).Value = "07:57:33"
    $ws.Cells.Item(5, 4).Value = "MarketMaking"
    $ws.Cells.Item(5, 5).Value = "DOWN"
    $ws.Cells.Item(5, 6).Value = 0.82
    $ws.Cells.Item(5, 7).Value = 0.8
    $ws.Cells.Item(5, 8).Value = "CLOSED"
    $ws.Cells.Item(5, 9).Value = -2.439
    $ws.Cells.Item(5, 10).Value = -0.02
    $ws.Cells.Item(5, 11).Value = 99.95
    $ws.Cells.Item(5, 12).Value = 0
    $ws.Cells.Item(5, 13).Value = 0
    $ws.Cells.Item(5, 14).Value = 0.6
    $ws.Cells.Item(5, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item(5, 16).Value = "early_exit"
    $ws.Cells.Item(5, 17).Value = 0.14
}

$wsAllTrades = $wb.Worksheets.Item("All Trades")
Add-TradeRow $wsAllTrades

$wsMarketMaking = $wb.Worksheets.Item("MarketMaking")
Add-TradeRow $wsMarketMaking
